$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs)
$ws.Range("G2").Value = 0.01490866666666667
$ws.Range("H2").Value = 0.044726
$ws.Range("M2").Value = 2.133443333333334
$ws.Range("N2").Value = 6.40033
$ws.Range("O2").Value = 0.2605947899689859
$ws.Range("P2").Value = 0.2605947899689859
$ws.Range("Q2").Value = 0.0318067955088889
$ws.Range("R2").Value = 0.28626115958
$ws.Range("S2").Value = 0.2605947899689859
$ws.Range("T2").Value = 0.2605947899689859

# Row 3 (Target cluster: FAPs)
$ws.Range("G3").Value = 0.01490866666666667
$ws.Range("H3").Value = 0.044726
$ws.Range("O3").Value = 0.5209338844846115
$ws.Range("P3").Value = 0.5209338844846116
$ws.Range("Q3").Value = 0.0635823822088889
$ws.Range("R3").Value = 0.5722414398800001
$ws.Range("S3").Value = 0.5209338844846115
$ws.Range("T3").Value = 0.5209338844846116

# Row 4 (Target cluster: MuSCs)
$ws.Range("G4").Value = 0.01490866666666667
$ws.Range("H4").Value = 0.044726
$ws.Range("M4").Value = 1.788586
$ws.Range("N4").Value = 5.365758
$ws.Range("O4").Value = 0.2184713255464024
$ws.Range("P4").Value = 0.2184713255464024
$ws.Range("Q4").Value = 0.02666543247866666
$ws.Range("R4").Value = 0.239988892308
$ws.Range("S4").Value = 0.2184713255464024
$ws.Range("T4").Value = 0.2184713255464024
